$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The word list is being trimmed: "romans" and "kocs" are dropped as target
# words, leaving only coffee / ants / otzi. Each of the three condition
# blocks (US / NZ / TR) had 5 rows (coffee, ants, romans, otzi, kocs); they
# shrink to 3 rows (coffee, ants, otzi).
#
# Before layout:
#   rows 2-6   -> US:  coffee, ants, romans, otzi, kocs
#   rows 7-11  -> NZ:  coffee, ants, romans, otzi, kocs
#   rows 12-16 -> TR:  coffee, ants, romans, otzi, kocs
#
# Delete the "romans" and "kocs" row from each block, working from the
# bottom of the sheet upward so the row numbers of rows still queued for
# deletion don't shift underneath us.
$ws.Rows("16").Delete()   # TR kocs
$ws.Rows("14").Delete()   # TR romans
$ws.Rows("11").Delete()   # NZ kocs
$ws.Rows("9").Delete()    # NZ romans
$ws.Rows("6").Delete()    # US kocs
$ws.Rows("4").Delete()    # US romans

# Selection moved from F19 to D20.
$ws.Range("D20").Select()
